# "hapus menu tidak terpakai" - add new "Brand" filter rows (Canon, FujiFilm,
# Nikon, Olympus, Panasonic, Pentax, Samsung, Sony, All Brand) right after the
# existing "Kamera Waterproof" row, highlighted with a red fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new blank rows at row 9 - shifts existing rows 9:30 down to 18:39.
$ws.Rows("9:17").Insert()

# Highlight the new rows with the red fill used to flag them.
$ws.Range("C9:F17").Interior.Color = 192

# Fill in the new brand rows (id in C, brand name in E).
$ws.Range("C9").Value = 274
$ws.Range("E9").Value = "Canon"

$ws.Range("C10").Value = 378
$ws.Range("E10").Value = "FujiFilm"

$ws.Range("C11").Value = 351
$ws.Range("E11").Value = "Nikon"

$ws.Range("C12").Value = 379
$ws.Range("E12").Value = "Olympus"

$ws.Range("C13").Value = 380
$ws.Range("E13").Value = "Panasonic "

$ws.Range("C14").Value = 381
$ws.Range("E14").Value = "Pentax"

$ws.Range("C15").Value = 382
$ws.Range("E15").Value = "Samsung"

$ws.Range("C16").Value = 383
$ws.Range("E16").Value = "Sony"

$ws.Range("C17").Value = 683
$ws.Range("E17").Value = "All Brand"

# Move the active selection to match the author's final selection.
$ws.Range("G18").Select()
